$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.404.94"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "2.100.78"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'334.14"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.5212"
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("D8").Value = "'0.4534"
$ws.Range("E8").Value = "  +3.81%  "
$ws.Range("D9").Value = "'54.53"
$ws.Range("E9").Value = "  +15.28%  "
$ws.Range("D10").Value = "'0.08892"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "'24.07"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").Value = "2.099.45"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "'6.797"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("E15").Value = "  +3.25%  "
$ws.Range("D16").Value = "'96.98"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "'0.00001143"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "'1.005"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'0.06619"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "'19.18"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'6.279"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").Value = "30.459.32"
$ws.Range("D24").Value = "'12.33"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").Value = "'2.338"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").Value = "2.341.10"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Value = "'22.18"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").Value = "'162.85"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("E29").Value = "  -3.70%  "
$ws.Range("D30").Value = "'133.03"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "'1.205"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "'0.1066"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").Value = "'6.365"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").Value = "'3.944"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("D37").Value = "'5.791"
$ws.Range("E37").Value = "  +5.34%  "
$ws.Range("D38").Value = "'0.02572"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").Value = "'0.06840"
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("D40").Value = "'0.2301"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "'12.71"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "'0.6865"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "'1.245"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").Value = "'2.316"
$ws.Range("E44").Value = "  +4.84%  "
$ws.Range("D45").Value = "'13.95"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").Value = "'0.6346"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").Value = "'3.652"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").Value = "'1.246"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "'0.00000000347"
$ws.Range("E49").Value = "  +17.90%  "
$ws.Range("D50").Value = "'83.01"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'1.201"
$ws.Range("E51").Value = "  +0.35%  "
